# testrapport.xlsx edit script
# - adds colour-coded "works / partially works / doesn't work" cells on the
#   Responsivitet sheet together with a small legend,
# - adds an "admin commands" row to the Funktionalitet sheet,
# - tweaks a column width and some view/print settings,
# - removes an old/garbage look (handled implicitly, nothing to delete here).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Responsivitet")
$ws2 = $wb.Worksheets.Item("Funktionalitet")

# ---------------------------------------------------------------------------
# Helper: copy a known style (green / red / yellow / wrap-text) from sheet2
# onto a destination cell without touching its value.
# ---------------------------------------------------------------------------
$greenSrc  = $ws2.Range("B2")   # fungerar (theme green fill)
$redSrc    = $ws2.Range("D2")   # fungerar inte (red fill)
$yellowSrc = $ws2.Range("C2")   # fungerar till viss del (yellow fill)
$wrapSrc   = $ws2.Range("A5")   # wrap-text, no fill

# ---------------------------------------------------------------------------
# Sheet1 "Responsivitet": colour the resolution/feature grid B2:E8
# Row legend: B = 360p, C = 780p, D = 1080p, E = 1920p
# ---------------------------------------------------------------------------
$redSrc.Copy()
$ws1.Range("B3:B8").PasteSpecial(-4122)

$yellowSrc.Copy()
$ws1.Range("C3:C8").PasteSpecial(-4122)

# New green fill (FF00B050) - first usage defines style index 5
$ws1.Range("B2").Interior.Color = 5287936
$ws1.Range("C2:D2").Interior.Color = 5287936
$ws1.Range("D3:D8").Interior.Color = 5287936
$ws1.Range("E3:E5").Interior.Color = 5287936
$ws1.Range("E7:E8").Interior.Color = 5287936

$yellowSrc.Copy()
$ws1.Range("E2").PasteSpecial(-4122)
$ws1.Range("E6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet1: legend explaining the colour coding (rows 11-13)
# ---------------------------------------------------------------------------
$ws1.Range("A11").Interior.Color = 5287936
$ws1.Range("B11").Value = "fungerar"

$yellowSrc.Copy()
$ws1.Range("A12").PasteSpecial(-4122)

$wrapSrc.Copy()
$ws1.Range("B12").PasteSpecial(-4122)
$ws1.Range("B12").Value = "fungerar till viss del"
$ws1.Range("B12").RowHeight = 43.5

$redSrc.Copy()
$ws1.Range("A13").PasteSpecial(-4122)

$wrapSrc.Copy()
$ws1.Range("B13").PasteSpecial(-4122)
$ws1.Range("B13").Value = "fungerar inte"
$ws1.Range("B13").RowHeight = 29

$excel.CutCopyMode = 0

$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Sheet2 "Funktionalitet": add a new "admin commands" row, widen column A
# ---------------------------------------------------------------------------
$ws2.Range("A9").Value = "admin commands"

$ws2.Range("B2").Copy()
$ws2.Range("B9").PasteSpecial(-4122)

$yellowSrc.Copy()
$ws2.Range("C9").PasteSpecial(-4122)
$ws2.Range("C9").Value = "x"

$ws2.Range("D2").Copy()
$ws2.Range("D9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws2.Columns.Item(1).ColumnWidth = 14.67
$ws2.Range("E9").Select()

# ---------------------------------------------------------------------------
# Re-activate sheet1 and restore its selection so it remains the tab shown
# when the workbook is opened (matches the original tabSelected state).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("E12").Select()
